$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we are about to write so that numeric-looking
# strings (e.g. "1.000", "13.06") are preserved as text, matching the source data,
# instead of being auto-converted into numbers by Excel. We reset the style back to
# "Normal" immediately after so no extra cell formatting/style is introduced.
$cells = @('D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E10', 'D11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'D18', 'E18', 'B19', 'C19', 'D19', 'E19', 'B20', 'C20', 'D20', 'E20', 'E21', 'E22', 'D23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'D28', 'E28', 'D29', 'E29', 'E30', 'D31', 'E31', 'D32', 'E32', 'D33', 'E33', 'D34', 'E34', 'E35', 'D36', 'E36', 'D37', 'E37', 'D38', 'E38', 'D39', 'E39', 'D40', 'E40', 'D41', 'E41', 'B42', 'C42', 'D42', 'E42', 'B43', 'C43', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45', 'E46', 'D47', 'E47', 'D48', 'E48', 'D49', 'E49', 'B50', 'C50', 'D50', 'E50', 'B51', 'C51', 'D51', 'E51')
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '29.131.87'
$ws.Range('E2').Value = '  -3.22%  '
$ws.Range('D3').Value = '1.849.88'
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').Value = '0.7042'
$ws.Range('E5').Value = '  -4.71%  '
$ws.Range('D6').Value = '238.11'
$ws.Range('E6').Value = '  -1.86%  '
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').Value = '0.3040'
$ws.Range('E8').Value = '  -4.00%  '
$ws.Range('D9').Value = '0.07509'
$ws.Range('E9').Value = '  +4.24%  '
$ws.Range('D10').Value = '23.38'
$ws.Range('E10').Value = '  -6.08%  '
$ws.Range('D11').Value = '0.08127'
$ws.Range('D12').Value = '0.7258'
$ws.Range('E12').Value = '  -4.44%  '
$ws.Range('D13').Value = '1.831.01'
$ws.Range('E13').Value = '  -5.53%  '
$ws.Range('D14').Value = '5.214'
$ws.Range('E14').Value = '  -4.13%  '
$ws.Range('D15').Value = '88.95'
$ws.Range('E15').Value = '  -4.14%  '
$ws.Range('D16').Value = '29.042.24'
$ws.Range('E16').Value = '  -3.73%  '
$ws.Range('D17').Value = '5.774'
$ws.Range('E17').Value = '  -6.62%  '
$ws.Range('D18').Value = '238.50'
$ws.Range('E18').Value = '  -4.73%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').Value = '13.06'
$ws.Range('E19').Value = '  -4.24%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.000007659'
$ws.Range('E20').Value = '  -2.62%  '
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('D23').Value = '2.074.72'
$ws.Range('D24').Value = '7.544'
$ws.Range('E24').Value = '  -5.65%  '
$ws.Range('D25').Value = '162.03'
$ws.Range('E25').Value = '  -1.36%  '
$ws.Range('D26').Value = '8.983'
$ws.Range('E26').Value = '  -3.38%  '
$ws.Range('D27').Value = '0.1458'
$ws.Range('E27').Value = '  -7.80%  '
$ws.Range('D28').Value = '18.02'
$ws.Range('E28').Value = '  -3.87%  '
$ws.Range('D29').Value = '1.940'
$ws.Range('E29').Value = '  -6.01%  '
$ws.Range('E30').Value = '  -6.38%  '
$ws.Range('D31').Value = '4.544'
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('D32').Value = '1.493'
$ws.Range('E32').Value = '  -2.80%  '
$ws.Range('D33').Value = '3.989'
$ws.Range('E33').Value = '  -5.38%  '
$ws.Range('D34').Value = '0.05146'
$ws.Range('E34').Value = '  -4.65%  '
$ws.Range('E35').Value = '  -5.21%  '
$ws.Range('D36').Value = '1.036'
$ws.Range('E36').Value = '  +3.47%  '
$ws.Range('D37').Value = '0.6997'
$ws.Range('E37').Value = '  -9.21%  '
$ws.Range('D38').Value = '2.644'
$ws.Range('E38').Value = '  -2.96%  '
$ws.Range('D39').Value = '0.01873'
$ws.Range('E39').Value = '  -4.92%  '
$ws.Range('D40').Value = '2.678'
$ws.Range('E40').Value = '  -3.17%  '
$ws.Range('D41').Value = '0.9443'
$ws.Range('E41').Value = '  +8.58%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '5.991'
$ws.Range('E42').Value = '  -1.42%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.080.20'
$ws.Range('E43').Value = '  -2.13%  '
$ws.Range('D44').Value = '0.4290'
$ws.Range('E44').Value = '  -5.86%  '
$ws.Range('D45').Value = '69.83'
$ws.Range('E45').Value = '  -3.84%  '
$ws.Range('E46').Value = '  -0.21%  '
$ws.Range('D47').Value = '102.32'
$ws.Range('E47').Value = '  -2.00%  '
$ws.Range('D48').Value = '1.743'
$ws.Range('E48').Value = '  -6.60%  '
$ws.Range('D49').Value = '1.980.76'
$ws.Range('E49').Value = '  -4.32%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '9.161'
$ws.Range('E50').Value = '  -4.56%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').Value = '7.041'
$ws.Range('E51').Value = '  -7.33%  '

foreach ($addr in $cells) {
    $ws.Range($addr).Style = "Normal"
}

